$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.284.81"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.610.46"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'212.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'18.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.836.47"
$ws.Range("D13").Value = "1.612.25"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "'4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "26.279.10"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "'62.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'201.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "'4.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").Value = "'9.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").Value = "'143.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "'15.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "'0.0500"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.62%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'3.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").Value = "'2.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "1.159.53"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("D37").Value = "'0.0166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D40").Value = "'0.789"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("E42").Value = "  +4.10%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "1.747.22"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "'92.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E46").Value = "  +13.34%  "
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "'53.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.27%  "
